# Update gh-pages generated output: bump "想去人数" (want-to-go) counts
# for a few events in both the "展览" sheet and the aggregated
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2, 4, 5
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 91
$wsExpo.Range("F4").Value = 4845
$wsExpo.Range("F5").Value = 16

# Sheet "全部类型" (all types) - same events, rows 2, 4, 6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 91
$wsAll.Range("F4").Value = 4845
$wsAll.Range("F6").Value = 16
